# Update "北京-漫展信息.xlsx" to the latest scraped data.
# Sheet 1 ("展览") gains a new row (new event on 2024-06-02) which shifts the
# existing rows 19-38 down to 20-39, and a number of "want to go" counts (column F)
# are refreshed. Sheet 4 ("全部类型") shares the same events but keeps its own
# row order, so only its F values need to be refreshed (no insert there).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: 展览 (Exhibition)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# --- Refresh F (want-to-go count) values for rows that are NOT affected by the
#     upcoming row insertion (rows 2-18 keep their row numbers). ---
$ws1.Cells.Item(8,  6).Value2 = 2262
$ws1.Cells.Item(9,  6).Value2 = 1499
$ws1.Cells.Item(10, 6).Value2 = 42
$ws1.Cells.Item(11, 6).Value2 = 828
$ws1.Cells.Item(12, 6).Value2 = 92
$ws1.Cells.Item(13, 6).Value2 = 2598
$ws1.Cells.Item(15, 6).Value2 = 1436
$ws1.Cells.Item(16, 6).Value2 = 5990
$ws1.Cells.Item(18, 6).Value2 = 5478

# --- Insert the new row for the "2024-06-02" event at row 19, pushing the old
#     rows 19-38 down to 20-39. ---
$ws1.Rows.Item(19).Insert()

$ws1.Cells.Item(19, 1).Value2 = 18
# Column B holds a date-looking string ("2024-06-02"); force it to stay text so
# Excel does not reinterpret it as a date serial number.
$ws1.Cells.Item(19, 2).Value2 = "'2024-06-02"
$ws1.Cells.Item(19, 2).Style = "Normal"
$ws1.Cells.Item(19, 3).Value2 = "北京·第二届CDS知名声优门胁舞以专场见面会"
$ws1.Cells.Item(19, 4).Value2 = "黑庄户路8号 北京音乐产业园"
$ws1.Cells.Item(19, 5).Value2 = "2024.06.02 10:30-06.02 14:00"
$ws1.Cells.Item(19, 6).Value2 = 0
$ws1.Cells.Item(19, 7).Value2 = 238
$ws1.Cells.Item(19, 8).Value2 = "https://show.bilibili.com/platform/detail.html?id=85389"
$ws1.Cells.Item(19, 9).Value2 = "//i1.hdslb.com/bfs/openplatform/202405/0GpdFuoU1715248227967.png"

# --- Refresh F values for the rows that shifted down (now at their new row
#     numbers, 20-39). ---
$ws1.Cells.Item(20, 6).Value2 = 2080
$ws1.Cells.Item(21, 6).Value2 = 2966
$ws1.Cells.Item(22, 6).Value2 = 3400
$ws1.Cells.Item(24, 6).Value2 = 1668
$ws1.Cells.Item(25, 6).Value2 = 36
$ws1.Cells.Item(26, 6).Value2 = 280
$ws1.Cells.Item(30, 6).Value2 = 344
$ws1.Cells.Item(31, 6).Value2 = 1062
$ws1.Cells.Item(32, 6).Value2 = 2248
$ws1.Cells.Item(36, 6).Value2 = 846

# ---------------------------------------------------------------------------
# Sheet 4: 全部类型 (All types) - same events, own row order, values only.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Cells.Item(12, 6).Value2 = 1499
$ws4.Cells.Item(13, 6).Value2 = 42
$ws4.Cells.Item(14, 6).Value2 = 828
$ws4.Cells.Item(15, 6).Value2 = 92
$ws4.Cells.Item(17, 6).Value2 = 2598
$ws4.Cells.Item(18, 6).Value2 = 1436
$ws4.Cells.Item(23, 6).Value2 = 5990
$ws4.Cells.Item(25, 6).Value2 = 5478
$ws4.Cells.Item(26, 6).Value2 = 2080
$ws4.Cells.Item(27, 6).Value2 = 2966
$ws4.Cells.Item(28, 6).Value2 = 3400
$ws4.Cells.Item(33, 6).Value2 = 1668
$ws4.Cells.Item(36, 6).Value2 = 280
$ws4.Cells.Item(40, 6).Value2 = 344
$ws4.Cells.Item(42, 6).Value2 = 2248
$ws4.Cells.Item(46, 6).Value2 = 846

$wb.Save()
